# Word COM-interop script implementing the commit
# "Se agrega César Colorado a Hoja de datos"
#
# Changes applied:
#   1. Add "César Alejandro Colorado Jacobo" / "Desarrollador" to the first
#      previously-empty row of the team table.
#   2. Split several names into separate runs bracketed by <w:proofErr/>
#      spell-check markers (as Word's spell checker does once it re-scans
#      text after an edit/open).
#   3. Move the `_GoBack` bookmark from the very last paragraph of the
#      document to right after "Modelo de Entidades " (Word drops this
#      bookmark at the location of the most-recent edit/cursor position).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Fill in the new team member row (previously blank cells).
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Item(8)

$nameCell = $newRow.Cells.Item(1).Range
$nameCell.End = $nameCell.End - 1
$nameCell.Text = "César Alejandro Colorado Jacobo"

$roleCell = $newRow.Cells.Item(2).Range
$roleCell.End = $roleCell.End - 1
$roleCell.Text = "Desarrollador"

# ---------------------------------------------------------------------
# 2) Split names into proofed runs (<w:proofErr w:type="spellStart"/> ...
#    <w:proofErr w:type="spellEnd"/>) around the word(s) Word's spell
#    checker doesn't recognise.
# ---------------------------------------------------------------------

function Set-ProofedRunsXml {
    param($SearchText, $InnerXml)
    $rng = $d.Content
    $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           $InnerXml +
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# "Denis Alfredo Vela Velasquez" -> "Denis Alfredo Vela " + proofed "Velasquez"
Set-ProofedRunsXml "Denis Alfredo Vela Velasquez" (
    '<w:pPr><w:spacing w:after="0"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Denis Alfredo Vela </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Velasquez</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# "Jose David Retana Retana" -> proofed "Jose" + " David Retana " + proofed "Retana"
Set-ProofedRunsXml "Jose David Retana Retana" (
    '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Jose</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> David Retana </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Retana</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# "Sebastian Eduardo Véliz Pinto" -> proofed "Sebastian" + " Eduardo Véliz Pinto" (es-ES lang)
Set-ProofedRunsXml "Sebastian Eduardo Véliz Pinto" (
    '<w:pPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Sebastian</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> Eduardo Véliz Pinto</w:t></w:r>'
)

# "Jose Luis Alejandro Estrada Hernández" -> proofed "Jose" + " Luis Alejandro Estrada Hernández"
Set-ProofedRunsXml "Jose Luis Alejandro Estrada Hernández" (
    '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Jose</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Luis Alejandro Estrada Hernández</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 3) Relocate the `_GoBack` bookmark.
# ---------------------------------------------------------------------

# Remove it from its old position (end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-insert it right after "Modelo de Entidades " via InsertXML (a
# collapsed Range.Add at this exact run-end boundary is unreliable, so we
# replace the trailing space character with itself plus the bookmark).
$spaceRng = $d.Content
$spaceRng.Find.Execute("Modelo de Entidades ", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0)
$spaceRng.Start = $spaceRng.End - 1
$bookmarkXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
               '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
               '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
               '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$spaceRng.InsertXML($bookmarkXml)

"Done"
